$wb = $excel.ActiveWorkbook

# ALC!row32
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H32").Value = 5469.6
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 5469.6
$ws.Range("K32").Value = 0
$ws.Range("L32").ClearContents()
$ws.Range("M32").ClearContents()
$ws.Range("N32").Value = -6121.6

# ALC!row47
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 12902.5
$ws.Range("I47").Value = 500
$ws.Range("J47").Value = 17036.666
$ws.Range("K47").Value = 500
$ws.Range("L47").Value = 17036.666
$ws.Range("M47").Value = 472
$ws.Range("N47").Value = -18980.666

# ALC!row52
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H52").Value = 1302.6666
$ws.Range("I52").Value = 512
$ws.Range("J52").Value = 10000
$ws.Range("K52").Value = 1536
$ws.Range("L52").Value = 30000
$ws.Range("M52").Value = -1376
$ws.Range("N52").Value = -30320

# ALC!row112
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1569.3448
$ws.Range("I112").Value = 485.6
$ws.Range("J112").Value = 1795.125
$ws.Range("K112").Value = 1456.8
$ws.Range("L112").Value = 5385.375
$ws.Range("M112").Value = -348.8000000000002
$ws.Range("N112").Value = -7601.375

# ALC!row137
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1538.2433
$ws.Range("I137").Value = 1618.1428
$ws.Range("J137").Value = 1289.6666
$ws.Range("K137").Value = 4854.428400000001
$ws.Range("L137").Value = 3868.9998
$ws.Range("M137").Value = -2304.428400000001
$ws.Range("N137").Value = -8968.9998

# ARM!row2
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1722.6097
$ws.Range("I2").Value = 1336.9231
$ws.Range("J2").Value = 2391.1333
$ws.Range("K2").Value = 1336.9231
$ws.Range("L2").Value = 2391.1333
$ws.Range("M2").Value = -1223.9231
$ws.Range("N2").Value = -2617.1333

# ARM!row45
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1660.037
$ws.Range("I45").Value = 1435.8235
$ws.Range("J45").Value = 2041.2
$ws.Range("K45").Value = 1435.8235
$ws.Range("L45").Value = 2041.2
$ws.Range("M45").Value = -1058.8235
$ws.Range("N45").Value = -2795.2

# ARM!row74
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 4264.5713
$ws.Range("I74").Value = 4094.4285
$ws.Range("J74").Value = 4349.643
$ws.Range("K74").Value = 4094.4285
$ws.Range("L74").Value = 4349.643
$ws.Range("M74").Value = -3220.4285
$ws.Range("N74").Value = -6097.643

# ARM!row77
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 4264.5713
$ws.Range("I77").Value = 4094.4285
$ws.Range("J77").Value = 4349.643
$ws.Range("K77").Value = 20472.1425
$ws.Range("L77").Value = 21748.215
$ws.Range("M77").Value = -16104.1425
$ws.Range("N77").Value = -30484.215

# ARM!row97
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H97").Value = 481.6842
$ws.Range("I97").Value = 333.18182
$ws.Range("J97").Value = 685.875
$ws.Range("K97").Value = 333.18182
$ws.Range("L97").Value = 685.875
$ws.Range("M97").Value = 162.81818
$ws.Range("N97").Value = -1677.875

# ARM!row116
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1722.6097
$ws.Range("I116").Value = 1336.9231
$ws.Range("J116").Value = 2391.1333
$ws.Range("K116").Value = 1336.9231
$ws.Range("L116").Value = 2391.1333
$ws.Range("M116").Value = 957.0769
$ws.Range("N116").Value = -6979.1333

# BSM!row3
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1722.6097
$ws.Range("I3").Value = 1336.9231
$ws.Range("J3").Value = 2391.1333
$ws.Range("K3").Value = 1336.9231
$ws.Range("L3").Value = 2391.1333
$ws.Range("M3").Value = -1222.9231
$ws.Range("N3").Value = -2619.1333

# CRP!row99
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 3421.932
$ws.Range("I99").Value = 2947.6177
$ws.Range("J99").Value = 5034.6
$ws.Range("K99").Value = 2947.6177
$ws.Range("L99").Value = 5034.6
$ws.Range("M99").Value = -1449.6177
$ws.Range("N99").Value = -8030.6

# CRP!row126
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 3421.932
$ws.Range("I126").Value = 2947.6177
$ws.Range("J126").Value = 5034.6
$ws.Range("K126").Value = 8842.8531
$ws.Range("L126").Value = 15103.8
$ws.Range("M126").Value = -6372.8531
$ws.Range("N126").Value = -20043.8

# CUL!row54
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H54").Value = 5666.3335
$ws.Range("I54").Value = 4000
$ws.Range("J54").Value = 6499.5
$ws.Range("K54").Value = 12000
$ws.Range("L54").Value = 19498.5
$ws.Range("M54").Value = -11441
$ws.Range("N54").Value = -20616.5

# CUL!row57
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H57").Value = 345335
$ws.Range("I57").Value = 345335
$ws.Range("J57").Value = 0
$ws.Range("K57").Value = 1036005
$ws.Range("L57").Value = 0
$ws.Range("M57").ClearContents()
$ws.Range("N57").ClearContents()

# CUL!row59
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H59").Value = 6508.909
$ws.Range("I59").Value = 1999
$ws.Range("J59").Value = 6959.9
$ws.Range("K59").Value = 5997
$ws.Range("L59").Value = 20879.7
$ws.Range("M59").Value = -5457
$ws.Range("N59").Value = -21959.7

# CUL!row62
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H62").Value = 1000
$ws.Range("I62").Value = 1000
$ws.Range("J62").Value = 0
$ws.Range("K62").Value = 3000
$ws.Range("L62").Value = 0
$ws.Range("M62").Value = -2314
$ws.Range("N62").ClearContents()

# CUL!row65
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H65").Value = 1000
$ws.Range("I65").Value = 1000
$ws.Range("J65").Value = 0
$ws.Range("K65").Value = 9000
$ws.Range("L65").Value = 0
$ws.Range("M65").Value = -5568
$ws.Range("N65").ClearContents()

# CUL!row69
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H69").Value = 6457.143
$ws.Range("I69").Value = 6457.143
$ws.Range("J69").Value = 0
$ws.Range("K69").Value = 19371.429
$ws.Range("L69").Value = 0
$ws.Range("M69").Value = -18560.429

# CUL!row72
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H72").Value = 6457.143
$ws.Range("I72").Value = 6457.143
$ws.Range("J72").Value = 0
$ws.Range("K72").Value = 58114.287
$ws.Range("L72").Value = 0
$ws.Range("M72").Value = -54058.287

# GSM!row102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2962.6155
$ws.Range("I102").Value = 2708.5557
$ws.Range("J102").Value = 3534.25
$ws.Range("K102").Value = 2708.5557
$ws.Range("L102").Value = 3534.25
$ws.Range("M102").Value = -1086.5557
$ws.Range("N102").Value = -6778.25

# GSM!row122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 3954.3076
$ws.Range("I122").Value = 2155.5557
$ws.Range("J122").Value = 8001.5
$ws.Range("K122").Value = 6466.6671
$ws.Range("L122").Value = 24004.5
$ws.Range("M122").Value = -4016.6671
$ws.Range("N122").Value = -28904.5

# GSM!row126
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4435.9287
$ws.Range("I126").Value = 4783.3335
$ws.Range("J126").Value = 4175.375
$ws.Range("K126").Value = 14350.0005
$ws.Range("L126").Value = 12526.125
$ws.Range("M126").Value = -11880.0005
$ws.Range("N126").Value = -17466.125

# GSM!row138
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H138").Value = 0
$ws.Range("I138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("L138").ClearContents()
$ws.Range("N138").ClearContents()

# LTW!row7
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5999.5884
$ws.Range("I7").Value = 3110.889
$ws.Range("J7").Value = 9249.375
$ws.Range("K7").Value = 3110.889
$ws.Range("L7").Value = 9249.375
$ws.Range("M7").Value = -2998.889
$ws.Range("N7").Value = -9473.375

# LTW!row16
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 548.9
$ws.Range("I16").Value = 548.9
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 548.9
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -378.9

# LTW!row82
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2356.7144
$ws.Range("I82").Value = 2356.7144
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 2356.7144
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -1995.7144
$ws.Range("N82").ClearContents()

# LTW!row85
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H85").Value = 2356.7144
$ws.Range("I85").Value = 2356.7144
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 2356.7144
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -1108.7144
$ws.Range("N85").ClearContents()

# LTW!row101
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H101").Value = 117832.4
$ws.Range("I101").Value = 0
$ws.Range("J101").Value = 117832.4
$ws.Range("K101").Value = 0
$ws.Range("L101").Value = 117832.4
$ws.Range("N101").Value = -124322.4

# LTW!row126
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H126").Value = 5999.5884
$ws.Range("I126").Value = 3110.889
$ws.Range("J126").Value = 9249.375
$ws.Range("K126").Value = 9332.667000000001
$ws.Range("L126").Value = 27748.125
$ws.Range("M126").Value = -6862.667000000001
$ws.Range("N126").Value = -32688.125

# LTW!row137
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H137").Value = 100000
$ws.Range("I137").Value = 0
$ws.Range("J137").Value = 100000
$ws.Range("K137").Value = 0
$ws.Range("L137").Value = 100000
$ws.Range("N137").Value = -110200

# WVR!row126
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H126").Value = 5960
$ws.Range("I126").Value = 6000
$ws.Range("J126").Value = 5900
$ws.Range("K126").Value = 18000
$ws.Range("L126").Value = 17700
$ws.Range("M126").Value = -15530
$ws.Range("N126").Value = -22640
